$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was M) -> B
$ws.Range("A2").Value = "B"
$ws.Range("B2").Value = 0.9452054794520548
$ws.Range("C2").Value = 0.965034965034965
$ws.Range("D2").Value = 0.9550173010380623
$ws.Range("E2").Value = 143

# Row 3 (was B) -> M
$ws.Range("A3").Value = "M"
$ws.Range("B3").Value = 0.9390243902439024
$ws.Range("C3").Value = 0.9058823529411765
$ws.Range("D3").Value = 0.9221556886227545
$ws.Range("E3").Value = 85

# Row 4 (accuracy)
$ws.Range("B4").Value = 0.9429824561403509
$ws.Range("C4").Value = 0.9429824561403509
$ws.Range("D4").Value = 0.9429824561403509
$ws.Range("E4").Value = 0.9429824561403509

# Row 5 (macro avg)
$ws.Range("B5").Value = 0.9421149348479786
$ws.Range("C5").Value = 0.9354586589880707
$ws.Range("D5").Value = 0.9385864948304083
$ws.Range("E5").Value = 228

# Row 6 (weighted avg)
$ws.Range("B6").Value = 0.942901126019191
$ws.Range("C6").Value = 0.9429824561403509
$ws.Range("D6").Value = 0.9427662613218292
$ws.Range("E6").Value = 228
